$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44894
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22500
$ws.Range("P4").Value = 22250
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("S4").Value = 2781
$ws.Range("T4").Value = 8

# Row 6
$ws.Range("D6").Value = 44533
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("S6").Value = 2312

# Row 7
$ws.Range("D7").Value = 44533

# Row 8
$ws.Range("D8").Value = 44498
$ws.Range("L8").Value = "Segunda"
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("S8").Value = 2438

# Row 9
$ws.Range("D9").Value = 44162
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 2000
$ws.Range("O9").Value = 2100
$ws.Range("P9").Value = 2050
$ws.Range("Q9").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S9").Value = 2050
$ws.Range("T9").Value = 1

# Row 14
$ws.Range("D14").Value = 44895
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 22500
$ws.Range("P14").Value = 22250
$ws.Range("Q14").Value = "$/bandeja 8 kilos"
$ws.Range("S14").Value = 2781
$ws.Range("T14").Value = 8

# Row 15
$ws.Range("D15").Value = 44526
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 21000
$ws.Range("S15").Value = 2625

# Row 16
$ws.Range("D16").Value = 44505
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("S16").Value = 2438

# Row 17
$ws.Range("D17").Value = 44488
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17500
$ws.Range("S17").Value = 2188

# Row 18
$ws.Range("D18").Value = 44530
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 19000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19500
$ws.Range("S18").Value = 2438

# Row 19
$ws.Range("D19").Value = 44530
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 16000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 16000
$ws.Range("S19").Value = 2000

# Row 20
$ws.Range("D20").Value = 44890
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 22500
$ws.Range("P20").Value = 22250
$ws.Range("S20").Value = 2781

# Row 21
$ws.Range("D21").Value = 44495
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 270
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19556
$ws.Range("S21").Value = 2444

# Row 23
$ws.Range("D23").Value = 44873
$ws.Range("M23").Value = 300

# Row 24
$ws.Range("D24").Value = 44491
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 18000
$ws.Range("O24").Value = 19000
$ws.Range("P24").Value = 18500
$ws.Range("S24").Value = 2312

# Row 25
$ws.Range("D25").Value = 44880
$ws.Range("L25").Value = "Primera"
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 22500
$ws.Range("P25").Value = 22250
$ws.Range("S25").Value = 2781

# Row 26
$ws.Range("D26").Value = 44159
$ws.Range("N26").Value = 2000
$ws.Range("O26").Value = 2100
$ws.Range("P26").Value = 2050
$ws.Range("Q26").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S26").Value = 2050
$ws.Range("T26").Value = 1
